$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 618014.4399999999
$ws.Range("I28").Value = 1111381.2
$ws.Range("J28").Value = 1306
$ws.Range("K28").Value = 1111381.2
$ws.Range("L28").Value = 1306
$ws.Range("M28").Value = -1110896.2
$ws.Range("N28").Value = -2276

$ws.Range("H111").Value = 1708.7778
$ws.Range("I111").Value = 1066.5
$ws.Range("J111").Value = 2993.3333
$ws.Range("K111").Value = 3199.5
$ws.Range("L111").Value = 8979.999899999999
$ws.Range("M111").Value = -132.5
$ws.Range("N111").Value = -15113.9999

$ws.Range("H115").Value = 1000
$ws.Range("I115").Value = 1000
$ws.Range("J115").Value = 0
$ws.Range("K115").Value = 3000
$ws.Range("L115").Value = 0
$ws.Range("M115").Value = -1433
$ws.Range("N115").ClearContents()

$ws.Range("H116").Value = 10646476
$ws.Range("I116").Value = 13839649
$ws.Range("J116").Value = 2566.6667
$ws.Range("K116").Value = 13839649
$ws.Range("L116").Value = 2566.6667
$ws.Range("M116").Value = -13836207
$ws.Range("N116").Value = -9450.6667

$ws.Range("H132").Value = 471891.88
$ws.Range("I132").Value = 578737.4
$ws.Range("J132").Value = 23140.8
$ws.Range("K132").Value = 1736212.2
$ws.Range("L132").Value = 69422.39999999999
$ws.Range("M132").Value = -1733682.2
$ws.Range("N132").Value = -74482.39999999999

$ws.Range("H133").Value = 42567.934
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 42567.934
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 42567.934
$ws.Range("N133").Value = -52687.934

$ws.Range("H134").Value = 62956.668
$ws.Range("I134").Value = 30000
$ws.Range("J134").Value = 69548
$ws.Range("K134").Value = 30000
$ws.Range("L134").Value = 69548
$ws.Range("M134").Value = -24930
$ws.Range("N134").Value = -79688

$ws.Range("H136").Value = 59393.332
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 59393.332
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 59393.332
$ws.Range("N136").Value = -69593.33199999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 32843.668
$ws.Range("I32").Value = 3978.6
$ws.Range("J32").Value = 177169
$ws.Range("K32").Value = 3978.6
$ws.Range("L32").Value = 177169
$ws.Range("M32").Value = -3691.6

$ws.Range("H110").Value = 2090.1052
$ws.Range("I110").Value = 1479.8
$ws.Range("J110").Value = 2308.0715
$ws.Range("K110").Value = 1479.8
$ws.Range("L110").Value = 2308.0715
$ws.Range("M110").Value = 565.2
$ws.Range("N110").Value = -6398.0715

$ws.Range("H139").Value = 0
$ws.Range("I139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H46").Value = 5300
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 5300
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 5300
$ws.Range("N46").Value = -5896

$ws.Range("H107").Value = 1401.625
$ws.Range("I107").Value = 1317.3889
$ws.Range("J107").Value = 1654.3334
$ws.Range("K107").Value = 1317.3889
$ws.Range("L107").Value = 1654.3334
$ws.Range("M107").Value = 602.6111000000001
$ws.Range("N107").Value = -5494.3334

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 801.3333
$ws.Range("I2").Value = 801.3333
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 801.3333
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -688.3333

$ws.Range("H16").Value = 938.5
$ws.Range("I16").Value = 941.3333
$ws.Range("J16").Value = 913
$ws.Range("K16").Value = 941.3333
$ws.Range("L16").Value = 913
$ws.Range("M16").Value = -654.3333

$ws.Range("H42").Value = 6355.5557
$ws.Range("I42").Value = 6000
$ws.Range("J42").Value = 6400
$ws.Range("K42").Value = 6000
$ws.Range("L42").Value = 6400
$ws.Range("M42").Value = -5407
$ws.Range("N42").Value = -7586

$ws.Range("H107").Value = 398.38095
$ws.Range("I107").Value = 226.22223
$ws.Range("J107").Value = 527.5
$ws.Range("K107").Value = 226.22223
$ws.Range("L107").Value = 527.5
$ws.Range("M107").Value = 1693.77777

$ws.Range("H113").Value = 938.5
$ws.Range("I113").Value = 941.3333
$ws.Range("J113").Value = 913
$ws.Range("K113").Value = 941.3333
$ws.Range("L113").Value = 913
$ws.Range("M113").Value = 1228.6667

$ws.Range("H132").Value = 2805.4119
$ws.Range("I132").Value = 2432.5518
$ws.Range("J132").Value = 4968
$ws.Range("K132").Value = 7297.655400000001
$ws.Range("L132").Value = 14904
$ws.Range("M132").Value = -4767.655400000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 2000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 2000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 6000
$ws.Range("N69").Value = -7622

$ws.Range("H72").Value = 2000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 2000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 18000
$ws.Range("N72").Value = -26112

$ws.Range("H114").Value = 1035
$ws.Range("I114").Value = 1694.6666
$ws.Range("J114").Value = 837.1
$ws.Range("K114").Value = 5083.9998
$ws.Range("L114").Value = 2511.3
$ws.Range("M114").Value = -1829.9998
$ws.Range("N114").Value = -9019.299999999999

$ws.Range("H121").Value = 482.2
$ws.Range("I121").Value = 196.66667
$ws.Range("J121").Value = 910.5
$ws.Range("K121").Value = 590.00001
$ws.Range("L121").Value = 2731.5
$ws.Range("M121").Value = 719.99999
$ws.Range("N121").Value = -5351.5

$ws.Range("H131").Value = 1354.9833
$ws.Range("I131").Value = 300
$ws.Range("J131").Value = 1372.8644
$ws.Range("K131").Value = 900
$ws.Range("L131").Value = 4118.593199999999
$ws.Range("M131").Value = 4140
$ws.Range("N131").Value = -14198.5932

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("M5").ClearContents()

$ws.Range("H107").Value = 915
$ws.Range("I107").Value = 1218
$ws.Range("J107").Value = 430.2
$ws.Range("K107").Value = 1218
$ws.Range("L107").Value = 430.2
$ws.Range("M107").Value = 702
$ws.Range("N107").Value = -4270.2

$ws.Range("H113").Value = 1759.6154
$ws.Range("I113").Value = 800
$ws.Range("J113").Value = 1934.091
$ws.Range("K113").Value = 800
$ws.Range("L113").Value = 1934.091
$ws.Range("M113").Value = 1370
$ws.Range("N113").Value = -6274.091

$ws.Range("H137").Value = 0
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H106").Value = 23611.285
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 23611.285
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 23611.285
$ws.Range("N106").Value = -26135.285

$ws.Range("H122").Value = 2948.0967
$ws.Range("I122").Value = 2156.5
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 6469.5
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -4019.5
$ws.Range("N122").Value = -15700

$ws.Range("H132").Value = 3653.2
$ws.Range("I132").Value = 2758.2
$ws.Range("J132").Value = 6338.2
$ws.Range("K132").Value = 8274.599999999999
$ws.Range("L132").Value = 19014.6
$ws.Range("M132").Value = -5744.599999999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H113").Value = 490
$ws.Range("I113").Value = 360
$ws.Range("J113").Value = 587.5
$ws.Range("K113").Value = 1080
$ws.Range("L113").Value = 1762.5
$ws.Range("M113").Value = 1090

$ws.Range("H122").Value = 84394.914
$ws.Range("I122").Value = 91794.45
$ws.Range("J122").Value = 3000
$ws.Range("K122").Value = 275383.35
$ws.Range("L122").Value = 9000
$ws.Range("M122").Value = -272933.35

$ws.Range("H133").Value = 56905
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 56905
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 56905
$ws.Range("N133").Value = -67025

$ws.Range("H136").Value = 2951.9
$ws.Range("I136").Value = 1523.1111
$ws.Range("J136").Value = 4120.909
$ws.Range("K136").Value = 4569.3333
$ws.Range("L136").Value = 12362.727
$ws.Range("M136").Value = -2019.3333
$ws.Range("N136").Value = -17462.727
